# JS and Py Minerva: Remove PaddingFrame and replace it with a CommentFrame.
#
# Row 6 of the "Protocol" sheet used to document PaddingFrame(numBytes); it
# now documents CommentFrame(comment) instead (same row/style, new text).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Protocol")

# --- A6: "CommentFrame(comment)" with "CommentFrame" bold, "(comment)" not ---
$ws.Range("A6").Value = "CommentFrame(comment)"

$boldPart = $ws.Range("A6").Characters(1, 12)
$boldPart.Font.Bold = $true
$boldPart.Font.Size = 10
$boldPart.Font.Name = "Tahoma"

$restPart = $ws.Range("A6").Characters(13, 9)
$restPart.Font.Bold = $false
$restPart.Font.Size = 10
$restPart.Font.Name = "Tahoma"

# --- F6: new description replacing the old padding-only description ---
$ws.Range("F6").Value = "CommentFrame is used for HTTP anti-script-inclusion preamble, padding, and heartbeats.  Padding is only needed to work around browser problems with content sniffing (in IE, Safari, Chrome, maybe Opera?), and maybe annoying proxies."

# --- Restore the view: scroll back to top-left and select F6 (matches the saved view state) ---
$ws.Activate()
$ws.Range("F6").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "applied CommentFrame edit"
